$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 43, shifting existing rows 43-128 down to 44-129
$ws.Rows.Item(43).Insert()

# Populate the newly inserted row 43 with the new weekly record
$ws.Cells.Item(43, 1).Value = 10
$ws.Cells.Item(43, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(43, 3).Value = "La Araucanía"
$ws.Cells.Item(43, 4).Value = 44645
$ws.Cells.Item(43, 5).Value = 9
$ws.Cells.Item(43, 6).Value = 100112012
$ws.Cells.Item(43, 7).Value = "Espinaca"
$ws.Cells.Item(43, 8).Value = "Sin especificar"
$ws.Cells.Item(43, 9).Value = "Primera"
$ws.Cells.Item(43, 10).Value = 65
$ws.Cells.Item(43, 11).Value = 10000
$ws.Cells.Item(43, 12).Value = 10000
$ws.Cells.Item(43, 13).Value = 10000
$ws.Cells.Item(43, 14).Value = "$/docena de atados"
$ws.Cells.Item(43, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(43, 16).Value = 3333
$ws.Cells.Item(43, 17).Value = 3
$ws.Cells.Item(43, 18).Value = "Hortaliza"
